$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# roboticS1Prep (column I) was recorded as text "No" for every sample row.
# Convert it to a real boolean (FALSE) with a TRUE/FALSE display format, so
# the column can be used/filtered like an actual yes-no flag going forward.
$boolFormat = '"TRUE";"TRUE";"FALSE"'

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 9)   # column I
    $cell.Value = $false
    $cell.NumberFormat = $boolFormat
}

# Reflect the edited column in the active selection, matching where the
# cursor was left after the cleanup pass.
$ws.Range("I2:I13").Select()
